$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows appended after the existing last row (226).
# Columns: A = date (serial), B = nuovi pos., C = somma mobile 7gg.,
#          D = somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @{ Row = 227; Date = 44301; B = 6;  C = 126; D = 174.1076980475065 },
    @{ Row = 228; Date = 44302; B = 26; C = 128; D = 176.8713123022289 },
    @{ Row = 229; Date = 44303; B = 14; C = 112; D = 154.7623982644503 }
)

$lastRow = 226

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Copy the formatting from the last existing row's date cell down into
    # the new row's date cell so it picks up the same cell style (date
    # style with border/bold/center alignment). Columns B/C/D keep the
    # default (unstyled) formatting, matching the workbook pattern.
    $ws.Range("A$lastRow").Copy($ws.Range("A$rowNum"))

    $ws.Cells.Item($rowNum, 1).Value = $r.Date
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D

    $lastRow = $rowNum
}
